$d = $word.ActiveDocument

# 1. Update the "Repositorio" link text to point to main branch.
$d.Content.Find.Execute(
    "https://github.com/HectorSanchezLuque/PIdam2122verde/",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "https://github.com/HectorSanchezLuque/PIdam2122verde/tree/main",
    2)

# 2. Remove the entire "Base de datos: ..." bullet paragraph (including
#    its paragraph mark), which the diff drops completely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Base de datos:*") {
        $p.Range.Delete()
        break
    }
}
